$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued cells (names, dismissal types, bowler names, overs-as-text) ---
$ws.Range('A2').NumberFormat = '@'
$ws.Range('A2').Value = 'Martin Guptill'
$ws.Range('A2').Style = 'Normal'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = 'Caught'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = ' Taskin Ahmed'
$ws.Range('E2').Style = 'Normal'
$ws.Range('J2').NumberFormat = '@'
$ws.Range('J2').Value = 'Tamim Iqbal'
$ws.Range('J2').Style = 'Normal'
$ws.Range('M2').NumberFormat = '@'
$ws.Range('M2').Value = 'Bowled'
$ws.Range('M2').Style = 'Normal'
$ws.Range('N2').NumberFormat = '@'
$ws.Range('N2').Value = ' Trent Boult'
$ws.Range('N2').Style = 'Normal'
$ws.Range('A3').NumberFormat = '@'
$ws.Range('A3').Value = 'Daryl Mitchell'
$ws.Range('A3').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = 'LBW'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = ' Mustafizur Rahman'
$ws.Range('E3').Style = 'Normal'
$ws.Range('J3').NumberFormat = '@'
$ws.Range('J3').Value = 'Liton Das'
$ws.Range('J3').Style = 'Normal'
$ws.Range('M3').NumberFormat = '@'
$ws.Range('M3').Value = 'Caught'
$ws.Range('M3').Style = 'Normal'
$ws.Range('N3').NumberFormat = '@'
$ws.Range('N3').Value = ' Mitchell Santner'
$ws.Range('N3').Style = 'Normal'
$ws.Range('A4').NumberFormat = '@'
$ws.Range('A4').Value = 'Kane Williamson(C)'
$ws.Range('A4').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = 'Bowled'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = ' Mahedi Hasan'
$ws.Range('E4').Style = 'Normal'
$ws.Range('J4').NumberFormat = '@'
$ws.Range('J4').Value = 'Shakib Al Hasan'
$ws.Range('J4').Style = 'Normal'
$ws.Range('M4').NumberFormat = '@'
$ws.Range('M4').Value = 'LBW'
$ws.Range('M4').Style = 'Normal'
$ws.Range('N4').NumberFormat = '@'
$ws.Range('N4').Value = ' Trent Boult'
$ws.Range('N4').Style = 'Normal'
$ws.Range('A5').NumberFormat = '@'
$ws.Range('A5').Value = 'Devon Conway'
$ws.Range('A5').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = 'LBW'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = ' Mahedi Hasan'
$ws.Range('E5').Style = 'Normal'
$ws.Range('J5').NumberFormat = '@'
$ws.Range('J5').Value = 'Mushfiqur Rahim'
$ws.Range('J5').Style = 'Normal'
$ws.Range('M5').NumberFormat = '@'
$ws.Range('M5').Value = 'LBW'
$ws.Range('M5').Style = 'Normal'
$ws.Range('N5').NumberFormat = '@'
$ws.Range('N5').Value = ' Tim Southee'
$ws.Range('N5').Style = 'Normal'
$ws.Range('A6').NumberFormat = '@'
$ws.Range('A6').Value = 'Glenn Phillips'
$ws.Range('A6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = ' Shamim Hossain'
$ws.Range('E6').Style = 'Normal'
$ws.Range('J6').NumberFormat = '@'
$ws.Range('J6').Value = 'Mahmudulla(C)'
$ws.Range('J6').Style = 'Normal'
$ws.Range('M6').NumberFormat = '@'
$ws.Range('M6').Value = 'Bowled'
$ws.Range('M6').Style = 'Normal'
$ws.Range('N6').NumberFormat = '@'
$ws.Range('N6').Value = ' Tim Southee'
$ws.Range('N6').Style = 'Normal'
$ws.Range('A7').NumberFormat = '@'
$ws.Range('A7').Value = 'James Neesham'
$ws.Range('A7').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = 'Caught'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = ' Shoriful Islam'
$ws.Range('E7').Style = 'Normal'
$ws.Range('J7').NumberFormat = '@'
$ws.Range('J7').Value = 'Afif Hossain'
$ws.Range('J7').Style = 'Normal'
$ws.Range('M7').NumberFormat = '@'
$ws.Range('M7').Value = 'Bowled'
$ws.Range('M7').Style = 'Normal'
$ws.Range('N7').NumberFormat = '@'
$ws.Range('N7').Value = ' Ish Sodhi'
$ws.Range('N7').Style = 'Normal'
$ws.Range('A8').NumberFormat = '@'
$ws.Range('A8').Value = 'Mitchell Santner'
$ws.Range('A8').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = 'Caught'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = ' Mustafizur Rahman'
$ws.Range('E8').Style = 'Normal'
$ws.Range('J8').NumberFormat = '@'
$ws.Range('J8').Value = 'Shamim Hossain'
$ws.Range('J8').Style = 'Normal'
$ws.Range('M8').NumberFormat = '@'
$ws.Range('M8').Value = 'Caught'
$ws.Range('M8').Style = 'Normal'
$ws.Range('N8').NumberFormat = '@'
$ws.Range('N8').Value = ' Mitchell Santner'
$ws.Range('N8').Style = 'Normal'
$ws.Range('A9').NumberFormat = '@'
$ws.Range('A9').Value = 'Adam Milne'
$ws.Range('A9').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = 'Bowled'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = ' Taskin Ahmed'
$ws.Range('E9').Style = 'Normal'
$ws.Range('J9').NumberFormat = '@'
$ws.Range('J9').Value = 'Mahedi Hasan'
$ws.Range('J9').Style = 'Normal'
$ws.Range('M9').NumberFormat = '@'
$ws.Range('M9').Value = 'Bowled'
$ws.Range('M9').Style = 'Normal'
$ws.Range('N9').NumberFormat = '@'
$ws.Range('N9').Value = ' Tim Southee'
$ws.Range('N9').Style = 'Normal'
$ws.Range('A10').NumberFormat = '@'
$ws.Range('A10').Value = 'Ish Sodhi'
$ws.Range('A10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = ' Taskin Ahmed'
$ws.Range('E10').Style = 'Normal'
$ws.Range('J10').NumberFormat = '@'
$ws.Range('J10').Value = 'Taskin Ahmed'
$ws.Range('J10').Style = 'Normal'
$ws.Range('M10').NumberFormat = '@'
$ws.Range('M10').Value = 'NOT OUT'
$ws.Range('M10').Style = 'Normal'
$ws.Range('A11').NumberFormat = '@'
$ws.Range('A11').Value = 'Tim Southee'
$ws.Range('A11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = ' Shoriful Islam'
$ws.Range('E11').Style = 'Normal'
$ws.Range('J11').NumberFormat = '@'
$ws.Range('J11').Value = 'Mustafizur Rahman'
$ws.Range('J11').Style = 'Normal'
$ws.Range('M11').NumberFormat = '@'
$ws.Range('M11').Value = 'LBW'
$ws.Range('M11').Style = 'Normal'
$ws.Range('N11').NumberFormat = '@'
$ws.Range('N11').Value = ' Mitchell Santner'
$ws.Range('N11').Style = 'Normal'
$ws.Range('A12').NumberFormat = '@'
$ws.Range('A12').Value = 'Trent Boult'
$ws.Range('A12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = 'NOT OUT'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = ' '
$ws.Range('E12').Style = 'Normal'
$ws.Range('J12').NumberFormat = '@'
$ws.Range('J12').Value = 'Shoriful Islam'
$ws.Range('J12').Style = 'Normal'
$ws.Range('M12').NumberFormat = '@'
$ws.Range('M12').Value = 'Caught'
$ws.Range('M12').Style = 'Normal'
$ws.Range('N12').NumberFormat = '@'
$ws.Range('N12').Value = ' Mitchell Santner'
$ws.Range('N12').Style = 'Normal'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = '11.0'
$ws.Range('C16').Style = 'Normal'
$ws.Range('L16').NumberFormat = '@'
$ws.Range('L16').Value = '9.4'
$ws.Range('L16').Style = 'Normal'
$ws.Range('A21').NumberFormat = '@'
$ws.Range('A21').Value = 'Mustafizur Rahman'
$ws.Range('A21').Style = 'Normal'
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = '2.0'
$ws.Range('B21').Style = 'Normal'
$ws.Range('J21').NumberFormat = '@'
$ws.Range('J21').Value = 'Trent Boult'
$ws.Range('J21').Style = 'Normal'
$ws.Range('K21').NumberFormat = '@'
$ws.Range('K21').Value = '2.0'
$ws.Range('K21').Style = 'Normal'
$ws.Range('A22').NumberFormat = '@'
$ws.Range('A22').Value = 'Taskin Ahmed'
$ws.Range('A22').Style = 'Normal'
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = '2.0'
$ws.Range('B22').Style = 'Normal'
$ws.Range('J22').NumberFormat = '@'
$ws.Range('J22').Value = 'Tim Southee'
$ws.Range('J22').Style = 'Normal'
$ws.Range('K22').NumberFormat = '@'
$ws.Range('K22').Value = '2.0'
$ws.Range('K22').Style = 'Normal'
$ws.Range('A23').NumberFormat = '@'
$ws.Range('A23').Value = 'Mahedi Hasan'
$ws.Range('A23').Style = 'Normal'
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = '2.0'
$ws.Range('B23').Style = 'Normal'
$ws.Range('J23').NumberFormat = '@'
$ws.Range('J23').Value = 'Ish Sodhi'
$ws.Range('J23').Style = 'Normal'
$ws.Range('K23').NumberFormat = '@'
$ws.Range('K23').Value = '2.0'
$ws.Range('K23').Style = 'Normal'
$ws.Range('A24').NumberFormat = '@'
$ws.Range('A24').Value = 'Shamim Hossain'
$ws.Range('A24').Style = 'Normal'
$ws.Range('J24').NumberFormat = '@'
$ws.Range('J24').Value = 'Adam Milne'
$ws.Range('J24').Style = 'Normal'
$ws.Range('K24').NumberFormat = '@'
$ws.Range('K24').Value = '2.0'
$ws.Range('K24').Style = 'Normal'
$ws.Range('A25').NumberFormat = '@'
$ws.Range('A25').Value = 'Shoriful Islam'
$ws.Range('A25').Style = 'Normal'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = '3.0'
$ws.Range('B25').Style = 'Normal'
$ws.Range('J25').NumberFormat = '@'
$ws.Range('J25').Value = 'Mitchell Santner'
$ws.Range('J25').Style = 'Normal'
$ws.Range('K25').NumberFormat = '@'
$ws.Range('K25').Value = '1.4'
$ws.Range('K25').Style = 'Normal'

# --- Numeric-valued cells (runs, balls, wickets, economy, etc.) ---
$ws.Range('B2').Value = 13
$ws.Range('K2').Value = 17
$ws.Range('L2').Value = 5
$ws.Range('B3').Value = 8
$ws.Range('C3').Value = 5
$ws.Range('K3').Value = 49
$ws.Range('L3').Value = 15
$ws.Range('B4').Value = 38
$ws.Range('C4').Value = 11
$ws.Range('K4').Value = 0
$ws.Range('L4').Value = 1
$ws.Range('B5').Value = 12
$ws.Range('C5').Value = 6
$ws.Range('K5').Value = 0
$ws.Range('L5').Value = 1
$ws.Range('B6').Value = 8
$ws.Range('K6').Value = 8
$ws.Range('L6').Value = 4
$ws.Range('B7').Value = 8
$ws.Range('K7').Value = 9
$ws.Range('L7').Value = 4
$ws.Range('B8').Value = 9
$ws.Range('C8').Value = 5
$ws.Range('K8').Value = 20
$ws.Range('L8').Value = 8
$ws.Range('B9').Value = 11
$ws.Range('C9').Value = 5
$ws.Range('K9').Value = 36
$ws.Range('L9').Value = 14
$ws.Range('B10').Value = 0
$ws.Range('C10').Value = 1
$ws.Range('K10').Value = 8
$ws.Range('L10').Value = 3
$ws.Range('B11').Value = 33
$ws.Range('C11').Value = 11
$ws.Range('L11').Value = 1
$ws.Range('B12').Value = 23
$ws.Range('C12').Value = 7
$ws.Range('L12').Value = 2
$ws.Range('A16').Value = 163
$ws.Range('D16').Value = 66
$ws.Range('J16').Value = 147
$ws.Range('K16').Value = 10
$ws.Range('M16').Value = 58
$ws.Range('C21').Value = 25
$ws.Range('D21').Value = 2
$ws.Range('E21').Value = 12.5
$ws.Range('L21').Value = 35
$ws.Range('M21').Value = 2
$ws.Range('N21').Value = 17.5
$ws.Range('C22').Value = 26
$ws.Range('D22').Value = 3
$ws.Range('E22').Value = 13
$ws.Range('L22').Value = 21
$ws.Range('M22').Value = 3
$ws.Range('N22').Value = 10.5
$ws.Range('C23').Value = 41
$ws.Range('D23').Value = 2
$ws.Range('E23').Value = 20.5
$ws.Range('L23').Value = 41
$ws.Range('N23').Value = 20.5
$ws.Range('C24').Value = 24
$ws.Range('D24').Value = 1
$ws.Range('E24').Value = 12
$ws.Range('L24').Value = 40
$ws.Range('M24').Value = 0
$ws.Range('N24').Value = 20
$ws.Range('C25').Value = 47
$ws.Range('E25').Value = 15.67
$ws.Range('M25').Value = 4
$ws.Range('N25').Value = 7.14
